# Update the cryptocurrency price list with the latest scraped values
# (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-23 and 27-51: update Price (D) and/or Volume(1h) (E) columns only.
# "Text" marks the Price values that look like plain numbers (e.g. "217.27")
# and therefore need the cell pre-formatted as Text so Excel's COM Value
# setter doesn't silently convert them to a numeric type - the source data
# keeps every Price cell as a text string (some even use "." as a thousands
# separator, e.g. "26.271.93", which is never at risk of conversion).
$updates = @(
    @{ Row = 2;  D = "26.271.93"; Text = $false; E = "  +1.85%  " },
    @{ Row = 3;  D = "1.648.75";  Text = $false; E = "  +0.80%  " },
    @{ Row = 4;  D = $null;       Text = $false; E = "  -0.49%  " },
    @{ Row = 5;  D = "217.27";    Text = $true;  E = "  +0.86%  " },
    @{ Row = 6;  D = "0.507";     Text = $true;  E = "  +0.69%  " },
    @{ Row = 7;  D = $null;       Text = $false; E = "  -0.54%  " },
    @{ Row = 8;  D = "0.259";     Text = $true;  E = "  -0.24%  " },
    @{ Row = 9;  D = "0.0638";    Text = $true;  E = "  -0.06%  " },
    @{ Row = 10; D = "20.16";     Text = $true;  E = "  +2.46%  " },
    @{ Row = 11; D = "0.0793";    Text = $true;  E = "  -0.14%  " },
    @{ Row = 12; D = "4.31";      Text = $true;  E = "  +0.78%  " },
    @{ Row = 13; D = "1.875.67";  Text = $false; E = "  +0.82%  " },
    @{ Row = 14; D = "1.634.98";  Text = $false; E = "  -0.04%  " },
    @{ Row = 15; D = "0.556";     Text = $true;  E = "  -0.97%  " },
    @{ Row = 16; D = "0.0₃0768";  Text = $false; E = "  -0.37%  " },
    @{ Row = 17; D = "63.76";     Text = $true;  E = "  +1.22%  " },
    @{ Row = 18; D = "26.252.83"; Text = $false; E = "  +1.76%  " },
    @{ Row = 19; D = "0.997";     Text = $true;  E = "  -0.45%  " },
    @{ Row = 20; D = $null;       Text = $false; E = "  +0.07%  " },
    @{ Row = 21; D = "194.73";    Text = $true;  E = "  +0.61%  " },
    @{ Row = 22; D = "10.09";     Text = $true;  E = "  +1.31%  " },
    @{ Row = 23; D = $null;       Text = $false; E = "  +1.02%  " },
    @{ Row = 27; D = "0.127";     Text = $true;  E = "  +2.95%  " },
    @{ Row = 28; D = "6.97";      Text = $true;  E = "  +0.67%  " },
    @{ Row = 29; D = "15.69";     Text = $true;  E = "  +0.72%  " },
    @{ Row = 30; D = "1.25";      Text = $true;  E = "  +1.07%  " },
    @{ Row = 31; D = "0.0504";    Text = $true;  E = "  +2.09%  " },
    @{ Row = 32; D = "3.37";      Text = $true;  E = "  +0.29%  " },
    @{ Row = 33; D = "3.27";      Text = $true;  E = "  +0.08%  " },
    @{ Row = 34; D = "1.61";      Text = $true;  E = "  +2.06%  " },
    @{ Row = 35; D = $null;       Text = $false; E = "  +0.65%  " },
    @{ Row = 36; D = "0.918";     Text = $true;  E = "  +1.70%  " },
    @{ Row = 37; D = "1.140.77";  Text = $false; E = "  +0.46%  " },
    @{ Row = 38; D = "0.556";     Text = $true;  E = "  +1.89%  " },
    @{ Row = 39; D = "2.51";      Text = $true;  E = "  -1.57%  " },
    @{ Row = 40; D = "0.0158";    Text = $true;  E = "  +1.50%  " },
    @{ Row = 41; D = $null;       Text = $false; E = "  -0.47%  " },
    @{ Row = 42; D = "5.65";      Text = $true;  E = "  +1.32%  " },
    @{ Row = 43; D = "100.62";    Text = $true;  E = "  -0.01%  " },
    @{ Row = 44; D = "0.801";     Text = $true;  E = "  -0.57%  " },
    @{ Row = 45; D = "1.783.56";  Text = $false; E = "  +0.74%  " },
    @{ Row = 46; D = "56.40";     Text = $true;  E = "  +2.14%  " },
    @{ Row = 47; D = "0.0₆0107";  Text = $false; E = "  -4.99%  " },
    @{ Row = 48; D = "1.50";      Text = $true;  E = "  +7.05%  " },
    @{ Row = 49; D = "0.0518";    Text = $true;  E = "  +2.94%  " },
    @{ Row = 50; D = "7.75";      Text = $true;  E = "  +3.32%  " },
    @{ Row = 51; D = "0.418";     Text = $true;  E = "  +0.30%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($u.Row, 4)
        if ($u.Text) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Rows 24-26: the coin ranking order changed, so the Coin/Link/Price/Volume
# values rotate between the three rows:
#   BinanceUSD, Toncoin, Monero  ->  Toncoin, Monero, BinanceUSD
$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.79"
$ws.Cells.Item(24, 5).Value = "  -1.44%  "

$ws.Cells.Item(25, 2).Value = "Monero"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "143.68"
$ws.Cells.Item(25, 5).Value = "  +0.38%  "

$ws.Cells.Item(26, 2).Value = "BinanceUSD"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "0.996"
$ws.Cells.Item(26, 5).Value = "  -0.59%  "
